# Apply the SO2_org sheet split + ferrous metal SO2 correction edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the SO2 sheet, place the copy right after it, and rename
#    the copy to "SO2_org" (this holds the original unmodified numbers).
# ---------------------------------------------------------------------
$so2 = $wb.Worksheets.Item("SO2")
$so2.Copy([System.Reflection.Missing]::Value, $so2)
$so2org = $wb.Worksheets.Item(2)
$so2org.Name = "SO2_org"

# ---------------------------------------------------------------------
# 2. Add the "ferrous metal SO2" extrapolated-fraction columns (K & L)
#    to the SO2_org sheet.
# ---------------------------------------------------------------------
# K4 picks up the same header-band formatting as the other row-4 labels
# (e.g. C4); L4 is a left-aligned variant of the same plain body font.
$so2org.Range("C4").Copy() | Out-Null
$so2org.Range("K4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$so2org.Range("K4").Value = "ferrous metal SO2"

$so2org.Range("L4").Font.Name = "Arial"
$so2org.Range("L4").Font.Size = 10
$so2org.Range("L4").HorizontalAlignment = -4131     # xlLeft
$so2org.Range("L4").VerticalAlignment = -4108       # xlCenter
$so2org.Range("L4").Value = "(extrapolated from previuos MEIC)"
$excel.CutCopyMode = 0

# K5:K12 / K13 simply take the worksheet's default (unstyled) look, which
# matches the plain numeric style already used as the default for columns
# outside the formatted A:H block.
$kValues = @{
    5  = 0.90540497991943358
    6  = 0.82913066026734716
    7  = 0.70299421270751949
    8  = 0.74642320026979103
    9  = 0.74878259329410557
    10 = 0.72535127128518384
    11 = 0.73674433853465293
    12 = 0.74864726127981651
}
foreach ($row in $kValues.Keys) {
    $so2org.Range("K$row").Value = $kValues[$row]
}
$so2org.Range("K13").Value = 0.80890148020000541

# ---------------------------------------------------------------------
# 3. Replace the "total" column D values on the SO2 sheet with formulas
#    that subtract the ferrous-metal SO2 fraction from the SO2_org
#    values, for rows 5-12.
# ---------------------------------------------------------------------
for ($row = 5; $row -le 12; $row++) {
    $so2.Range("D$row").Formula = "=SO2_org!D$row-SO2_org!K$row"
}

# ---------------------------------------------------------------------
# 4. Leave each sheet's selection the way the author ended up with,
#    then make SO2 the active/selected sheet again.
# ---------------------------------------------------------------------
$so2org.Select() | Out-Null
$so2org.Range("D10").Select() | Out-Null

$so2.Select() | Out-Null
$so2.Range("H16").Select() | Out-Null
